$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Msg8706")

# B and F columns are constant across all new rows, matching existing shared strings
$bVal = "8: 255`n"
$fVal = "10F872226797"

$data = @"
924|2022-10-24 13:43:45|3600|1666611832|46|202|16
925|2022-10-24 13:44:06|3600|1666611854|46|208|16
926|2022-10-24 13:45:06|3600|1666611913|46|188|16
927|2022-10-24 13:46:06|3600|1666611973|46|188|16
928|2022-10-24 13:47:06|3600|1666612033|46|189|16
929|2022-10-24 13:48:06|3600|1666612093|46|187|16
930|2022-10-24 13:49:06|3600|1666612153|46|183|16
931|2022-10-27 22:39:04|3600|1666903154|50|-58756|3
932|2022-10-27 22:39:08|3600|1666903155|50|280|3
933|2022-10-27 22:40:05|3600|1666903215|50|336|3
934|2022-10-27 22:41:05|3600|1666903275|50|248|3
935|2022-10-27 22:42:05|3600|1666903335|50|561|3
936|2022-10-27 22:43:05|3600|1666903395|50|396|3
937|2022-10-27 22:44:06|3600|1666903456|50|197|3
938|2022-10-27 22:45:05|3600|1666903516|50|177|3
939|2022-10-27 22:46:05|3600|1666903576|50|195|3
940|2022-10-27 22:47:05|3600|1666903636|50|219|3
941|2022-10-27 22:48:05|3600|1666903696|50|234|3
942|2022-10-27 22:49:05|3600|1666903756|50|467|3
943|2022-10-27 22:50:05|3600|1666903816|50|182|3
944|2022-10-27 22:51:05|3600|1666903876|50|247|3
945|2022-10-27 22:52:05|3600|1666903936|50|421|3
946|2022-10-27 22:53:05|3600|1666903996|50|196|3
947|2022-10-27 22:54:05|3600|1666904056|50|185|3
948|2022-10-27 22:55:05|3600|1666904116|50|193|3
949|2022-10-27 22:56:06|3600|1666904176|50|192|3
950|2022-10-27 22:57:06|3600|1666904236|50|1296|3
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "\|"
    $r = [int]$parts[0]
    $aVal = $parts[1]
    $cVal = [double]$parts[2]
    $dVal = [double]$parts[3]
    $eVal = [double]$parts[4]
    $gVal = [double]$parts[5]
    $hVal = [double]$parts[6]

    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
    $ws.Cells.Item($r, 6).Value = $fVal
    $ws.Cells.Item($r, 7).Value = $gVal
    $ws.Cells.Item($r, 8).Value = $hVal
}

# The B column text contains an embedded newline; Excel would otherwise mark
# these rows with a custom wrapped-text row height. Re-running AutoFit
# restores the default (non-custom) row height, matching the source data.
$ws.Range("A924:H950").Rows.AutoFit()

Write-Host "Done: added rows 924-950"
